$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.316.90"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.571.20"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'585.12"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").Value = "'148.33"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  +4.33%  "
$ws.Range("D10").Value = "'5.69"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "'27.57"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "3.031.47"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "63.301.53"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("D17").Value = "2.575.77"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "'11.38"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'343.32"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").Value = "'4.42"
$ws.Range("E20").Value = "  +3.61%  "
$ws.Range("D21").Value = "'6.87"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'66.75"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "2.687.05"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'1.65"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'8.24"
$ws.Range("E27").Value = "  +14.60%  "
$ws.Range("D28").Value = "'8.55"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  +8.86%  "
$ws.Range("D32").Value = "0.0₃0828"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "'457.65"
$ws.Range("E33").Value = "  +13.44%  "
$ws.Range("D34").Value = "'1.64"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'176.77"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.408"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'19.29"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").Value = "'4.52"
$ws.Range("E38").Value = "  +4.89%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'1.75"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'151.87"
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("D44").Value = "'21.20"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  +7.50%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "'18.50"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'11.39"
$ws.Range("E51").Value = "  -0.06%  "
